# Apply edits to the "geometry calculation" sheet of the 2-storey generic
# home geometry workbook, adding new retrofit materials data that changes
# the footprint area and number of floors inputs, plus updated wall
# constant used in the H15/H18 formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("geometry calculation")

# Update the footprint area input (C2): 285.3 -> 264.3
$ws.Range("C2").Value = 264.3

# Update the number of floors input (C3): 3 -> 2
$ws.Range("C3").Value = 2

# Update the constant term used in the H15 / H18 wall-length formulas
# (10.363 -> 7.9248), keeping the rest of the formula intact.
# NOTE: single-quoted strings are used so that PowerShell does not try to
# expand "$C" / "$5" as variables.
$ws.Range("H15").Formula = '=7.9248+($C$5*H14)'
$ws.Range("H18").Formula = '=7.9248+($C$5*H17)'

# Move the active selection to C18, as recorded in the saved view state.
$ws.Range("C18").Select()

$excel.Calculate()
